$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6-9: LGBM-2500-20-37 block ---
$ws.Range("A6").Value = "LGBM-2500-20-37"
$ws.Range("B6").Value = "SMOTE+RUS(0.1,0.5), 0.1 split"
$ws.Range("C6").Value = 0.85
$ws.Range("D6").Value = 0.9
$ws.Range("E6").Value = 0.88
$ws.Range("F6").Value = 0.92

$ws.Range("A7").Value = "LGBM-2500-20-37"
$ws.Range("B7").Value = "SMOTE+RUS, 0.1 split (testing on original data) (training 62k rows each)"
$ws.Range("C7").Value = 0.63
$ws.Range("D7").Value = 0.72
$ws.Range("E7").Value = 0.67
$ws.Range("F7").Value = 0.77
$ws.Rows.Item(7).RowHeight = 15

$ws.Range("A8").Value = "LGBM-2500-20-37"
$ws.Range("B8").Value = "SMOTE+RUS, 0.1 split (testing on original data) (training 90k rows each)"
$ws.Range("C8").Value = 0.57999999999999996
$ws.Range("D8").Value = 0.68
$ws.Range("E8").Value = 0.63
$ws.Range("F8").Value = 0.8

$ws.Range("A9").Value = "LGBM-2500-20-37"
$ws.Range("B9").Value = "SMOTE+RUS, 0.1 split (testing on original data) (training 72k rows each)"
$ws.Range("C9").Value = 0.62
$ws.Range("D9").Value = 0.72
$ws.Range("E9").Value = 0.67
$ws.Range("F9").Value = 0.78

# --- Row 14: RF-100-5 block ---
$ws.Range("A14").Value = "RF-100-5"
$ws.Range("B14").Value = "original data; 0.25 split"
$ws.Range("C14").Value = 0.56999999999999995
$ws.Range("D14").Value = 0.62
$ws.Range("E14").Value = 0.6

# --- Row 15-19: RF-1500-20 block ---
$ws.Range("A15").Value = "RF-1500-20"
$ws.Range("B15").Value = "original data; 0.1 split; 100% balanced"
$ws.Range("C15").Value = 0.78
$ws.Range("D15").Value = 0.72
$ws.Range("E15").Value = 0.75
$ws.Range("F15").Value = 0.74

$ws.Range("A16").Value = "RF-1500-20"
$ws.Range("B16").Value = "original data; 0.25 split; 100% balanced"
$ws.Range("C16").Value = 0.77
$ws.Range("D16").Value = 0.72
$ws.Range("E16").Value = 0.74
$ws.Range("F16").Value = 0.73

$ws.Range("A17").Value = "RF-1500-20"
$ws.Range("B17").Value = "SMOTE+RUS, 0.1 split (testing on original data) (training 62k rows each)"
$ws.Range("C17").Value = 0.68
$ws.Range("D17").Value = 0.69
$ws.Range("E17").Value = 0.68
$ws.Range("F17").Value = 0.76

$ws.Range("A18").Value = "RF-1500-20"
$ws.Range("B18").Value = "SMOTE+RUS, 0.1 split (testing on original data) (training 72k rows each)"
$ws.Range("C18").Value = 0.66
$ws.Range("D18").Value = 0.66
$ws.Range("E18").Value = 0.66
$ws.Range("F18").Value = 0.77

$ws.Range("A19").Value = "RF-1500-20"
$ws.Range("B19").Value = "SMOTE+RUS, 0.1 split (testing on original data) (training 90k rows each)"
$ws.Range("C19").Value = 0.63
$ws.Range("D19").Value = 0.63
$ws.Range("E19").Value = 0.63
$ws.Range("F19").Value = 0.78

# --- Row 22-24: BalancedRF-1500-20 block ---
$ws.Range("A22").Value = "BalancedRF-1500-20"
$ws.Range("B22").Value = "original data;0.1 split; 100k sample of class 0"
$ws.Range("C22").Value = 0.77
$ws.Range("D22").Value = 0.53
$ws.Range("E22").Value = 0.63
$ws.Range("F22").Value = 0.73

$ws.Range("A23").Value = "BalancedRF-1500-20"
$ws.Range("B23").Value = "original data;0.1 split; 200k sample of class 0"
$ws.Range("C23").Value = 0.77
$ws.Range("D23").Value = 0.38
$ws.Range("E23").Value = 0.5
$ws.Range("F23").Value = 0.73

$ws.Range("A24").Value = "BalancedRF-1500-20"
$ws.Range("B24").Value = "original data;0.1 split; 60k sample of class 0"
$ws.Range("C24").Value = 0.77
$ws.Range("D24").Value = 0.64
$ws.Range("E24").Value = 0.7
$ws.Range("F24").Value = 0.73

# --- Column B width widened to fit the longer strings ---
$ws.Columns.Item(2).ColumnWidth = 62.33

# --- Selection moved to B14, matching the final saved cursor position ---
$ws.Range("B14").Select()
